# Update the "grid_cell" (AG) column on the "solar" sheet so that each
# distr_solelc_won-CHE_00xx process row points at its new grid cell id.
# (Rows 24 and 25 are unchanged in the source data.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$gridCellByRow = @{
    4  = "CHE_10"
    5  = "CHE_22"
    6  = "CHE_0"
    7  = "CHE_7"
    8  = "CHE_20"
    9  = "CHE_1"
    10 = "CHE_6"
    11 = "CHE_24"
    12 = "CHE_8"
    13 = "CHE_5"
    14 = "CHE_13"
    15 = "CHE_12"
    16 = "CHE_2"
    17 = "CHE_3"
    18 = "CHE_9"
    19 = "CHE_21"
    20 = "CHE_4"
    21 = "CHE_17"
    22 = "CHE_19"
    23 = "CHE_23"
    26 = "CHE_11"
    27 = "CHE_15"
    28 = "CHE_25"
}

foreach ($row in $gridCellByRow.Keys) {
    $ws.Range("AG$row").Value = $gridCellByRow[$row]
}
